# Update column F ("想去人数" / want-to-go count) values across sheets
# to match the regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 272
$ws.Range("F4").Value = 498
$ws.Range("F5").Value = 2307
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 62
$ws.Range("F8").Value = 77
$ws.Range("F9").Value = 1677
$ws.Range("F10").Value = 1677
$ws.Range("F11").Value = 1384
$ws.Range("F12").Value = 74
$ws.Range("F16").Value = 715
$ws.Range("F17").Value = 180
$ws.Range("F19").Value = 7421
$ws.Range("F20").Value = 8293
$ws.Range("F28").Value = 269
$ws.Range("F34").Value = 1480
$ws.Range("F40").Value = 764
$ws.Range("F44").Value = 263
$ws.Range("F45").Value = 211
$ws.Range("F47").Value = 199
$ws.Range("F48").Value = 181

# --- Sheet "演出" (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 307

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 148
$ws.Range("F6").Value = 20

# --- Sheet "全部类型" (all types, aggregate of the above) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 148
$ws.Range("F7").Value = 498
$ws.Range("F8").Value = 2307
$ws.Range("F9").Value = 6
$ws.Range("F10").Value = 62
$ws.Range("F11").Value = 77
$ws.Range("F12").Value = 1677
$ws.Range("F13").Value = 1677
$ws.Range("F14").Value = 20
$ws.Range("F15").Value = 74
$ws.Range("F19").Value = 715
$ws.Range("F21").Value = 180
$ws.Range("F23").Value = 7421
$ws.Range("F24").Value = 7421
$ws.Range("F25").Value = 8293
$ws.Range("F29").Value = 269
$ws.Range("F41").Value = 764
$ws.Range("F46").Value = 263
$ws.Range("F47").Value = 211
$ws.Range("F48").Value = 199
$ws.Range("F49").Value = 181
$ws.Range("F50").Value = 307
